# Update "想去人数" (column F) counts for a batch of events, and rename
# one event, mirroring the same edits across the "展览" sheet and the
# "全部类型" aggregate sheet (which both list these same events).

$wb = $excel.ActiveWorkbook

# event name (as currently stored in column C) -> new F value
$updates = @{
    "常熟·ACG动漫游戏嘉年华"                               = 515
    "常熟·CDW.动漫展05"                                    = 1510
    "昆山·2024随机宅舞&正反派对决（免费活动）"             = 149
    "苏州·国乙主场·次元燃歌·偶像时刻（免费展）"           = 142
    "苏州·授渔仲夏动漫节2.0"                               = 732
    "昆山·2024首届华盟次元动漫嘉年华（免费活动）"         = 323
    "苏州·OrangeOrange夏日随舞派对【免费展会】"           = 50
    "苏州·第三届.OCG.Summer Carnival-国潮动漫游戏嘉年华"  = 6369
    "苏州·艾卡动漫游戏嘉年华（免票展）"                   = 3
    "苏州·OCG国潮动漫游戏嘉年华凌飞内场"                 = 148
    "苏州·ICAN summer World动漫品牌夏游节"                 = 15248
    "苏州·第二届Redamancy动漫游戏嘉年华"                  = 1511
    "苏州·排球少年only-茶歇"                               = 277
    "苏州·星梦X动漫游戏展（免费展）"                       = 137
    "苏州·Good jump ACG中秋嘉年华动漫国潮文化节"          = 11021
    "苏州·I COME ACG动漫品牌博览会"                        = 746
    "苏州·第十三届理想乡动漫展-同人创作者大会"             = 4310
    "苏州·第四届-OCG国朝动漫游戏嘉年华"                    = 232
}

$oldName = "苏州·第十三届理想乡动漫展-同人创作者大会"
$newName = "苏州·理想乡动漫游戏展-两馆全开+三馆间通道"

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($row = 2; $row -le $lastRow; $row++) {
        $name = $ws.Cells.Item($row, 3).Value()

        if ($updates.ContainsKey($name)) {
            $ws.Cells.Item($row, 6).Value = $updates[$name]
        }

        if ($name -eq $oldName) {
            $ws.Cells.Item($row, 3).Value = $newName
        }
    }
}
